# Automatic update of files.
#
# 1) Column C ("Förändrad") on every data row (2-117) moves from the old
#    serial date 45184 (2023-09-15) to the new serial date 45186
#    (2023-09-17).
# 2) Every HYPERLINK() formula in columns S:Y gets a second argument
#    (the "friendly name") equal to the case id in column A of that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldChanged = 45184
$newChanged = 45186

$firstRow = 2
$lastRow = 117
$hyperlinkCols = @("S", "T", "U", "V", "W", "X", "Y")

for ($row = $firstRow; $row -le $lastRow; $row++) {

    $cCell = $ws.Cells.Item($row, 3)   # column C
    if ($cCell.Value2 -eq $oldChanged) {
        $cCell.Value = $newChanged
    }

    $label = $ws.Cells.Item($row, 1).Value2   # column A, e.g. "A 13663-2023"

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range($col + $row)
        $formula = $cell.Formula
        if ([string]::IsNullOrEmpty($formula)) {
            continue
        }
        if ($formula.StartsWith("=HYPERLINK(") -and -not $formula.Contains(",")) {
            $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $label + '")'
            $cell.Formula = $newFormula
        }
    }
}
